# Apply correction to AV-MCPS (row 3) and MCPS (row 9) Diebold-Mariano test results
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: AV-MCPS
$ws.Range("B3").Value = 2.046924673881068
$ws.Range("C3").Value = 0.614457099736811
$ws.Range("D3").Value = 1.432467574144257
$ws.Range("E3").Value = 69.98145033975428
$ws.Range("F3").Value = 13473

# Row 9: MCPS
$ws.Range("B9").Value = 2.541979610502438
$ws.Range("C9").Value = 0.6121381791025303
$ws.Range("D9").Value = 1.929841431399907
$ws.Range("E9").Value = 75.9188399240725
$ws.Range("F9").Value = 14723

# G3 / G9 need scientific-notation magnitudes (e-92 / e-90) that overflow
# plain numeric literals in this parser, so set them via a formula string
# and then flatten the formula down to a static value with PasteSpecial.
$ws.Range("G3").Formula = "=1.550017749819525E-92"
$ws.Range("G3").Copy()
$ws.Range("G3").PasteSpecial(-4163)

$ws.Range("G9").Formula = "=1.790797395140314E-90"
$ws.Range("G9").Copy()
$ws.Range("G9").PasteSpecial(-4163)

$excel.CutCopyMode = 0
